$d = $word.ActiveDocument

# 1) Version number 1.0 -> 1.2.5
$d.Content.Find.Execute("1.0", $true, $false, $false, $false, $false, $true, 1, $false, "1.2.5", 2) | Out-Null

# 2) Creation -> Update
$d.Content.Find.Execute("Creation", $true, $false, $false, $false, $false, $true, 1, $false, "Update", 2) | Out-Null

# 3) Author name Fabrício Araújo -> Julio Paiva
$d.Content.Find.Execute("Fabrício Araújo", $true, $false, $false, $false, $false, $true, 1, $false, "Julio Paiva", 2) | Out-Null

# 4) Date 27/08/2020 -> 31/05/2023
$d.Content.Find.Execute("27/08/2020", $true, $false, $false, $false, $false, $true, 1, $false, "31/05/2023", 2) | Out-Null

# 5) "O usuario devidamente autenticado e na tela de listagem de diárias" -> "O usuário devidamente autenticado e na tela de listagem de diárias."
$d.Content.Find.Execute("O usuario devidamente autenticado e na tela de listagem de diárias", $true, $false, $false, $false, $false, $true, 1, $false, "O usuário devidamente autenticado e na tela de listagem de diárias.", 2) | Out-Null

# 6) extratificação -> estratificação
$d.Content.Find.Execute("extratificação", $true, $false, $false, $false, $false, $true, 1, $false, "estratificação", 2) | Out-Null

# 7) numero do empenho -> número do empenho
$d.Content.Find.Execute("Apresenta o numero do empenho", $true, $false, $false, $false, $false, $true, 1, $false, "Apresenta o número do empenho", 2) | Out-Null
